$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$s.Delete()
